{"js": "// Office.js (Word JavaScript API) script\n// Goal (per the diff):\n//   1. Move the \"PUMP:RISK:10..50\" block (5 paragraphs) from just before the\n//      \"PUMP:UNIT:100\" block to just after the intro paragraph (i.e. right\n//      before \"PUMP:URS:1\").\n//   2. Move the trailing block of 8 paragraphs\n//        PUMP:HRS:103, PUMP:PRS:103, ACE:SRS:110, ACE:SRS:120,\n//        PUMP:TBV:1, PUMP:PRS:6, PUMP:DER:2, ACE:SRS:1000\n//      (currently at the very end, after \"PUMP:UNIT:220\") to just after\n//      \"PUMP:URS:4000\" (i.e. right before the \"PUMP:UNIT:100\" block), and\n//      while doing so, insert two brand-new lines into that block:\n//        - \"PUMP:TBV:1111\" right after \"PUMP:HRS:103\"\n//        - \"PUMP:TBD:1\" right after \"PUMP:PRS:6\" (before \"PUMP:DER:2\")\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Helper: find the (first) paragraph whose text equals `t` exactly.\nfunction findByText(items, t) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === t) return items[i];\n  }\n  throw new Error(\"paragraph not found: \" + JSON.stringify(t));\n}\n\nconst items = paragraphs.items;\n\n// --- Step 1: insert the new RISK block right before \"PUMP:URS:1\" ---\n// (Inserting *before* the plain/non-bold \"PUMP:URS:1 \" paragraph - rather\n// than *after* the bold intro paragraph - means the new runs don't inherit\n// the intro's bold run formatting, matching the diff where these new\n// paragraphs carry no <w:rPr> at all.)\nconst urs1Para = findByText(items, \"PUMP:URS:1 \");\n\nconst riskTexts = [\n  \"PUMP:RISK:10 \",\n  \"PUMP:RISK:20 \",\n  \"PUMP:RISK:30 \",\n  \"PUMP:RISK:40 \",\n  \"PUMP:RISK:50 \",\n];\n\n// Keep the anchor fixed (always insert right before \"PUMP:URS:1\") so the\n// paragraphs end up in the same order as `riskTexts`.\nfor (const t of riskTexts) {\n  urs1Para.insertParagraph(t, Word.InsertLocation.before);\n}\nawait context.sync();\n\n// --- Step 2: insert the moved/expanded HRS..ACE block right after PUMP:URS:4000 ---\n// Re-load paragraphs since the collection changed.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nconst urs4000Para = findByText(paragraphs2.items, \"PUMP:URS:4000 \");\n\nconst movedBlockTexts = [\n  \"PUMP:HRS:103\",\n  \"PUMP:TBV:1111\", // new line\n  \"PUMP:PRS:103\",\n  \"ACE:SRS:110\",\n  \"ACE:SRS:120\",\n  \"PUMP:TBV:1\",\n  \"PUMP:PRS:6\",\n  \"PUMP:TBD:1\", // new line\n  \"PUMP:DER:2\",\n  \"ACE:SRS:1000\",\n];\n\nlet anchor = urs4000Para;\nfor (const t of movedBlockTexts) {\n  anchor = anchor.insertParagraph(t, Word.InsertLocation.after);\n}\nawait context.sync();\n\n// --- Step 3: delete the original RISK block (the one that was right before\n//     the PUMP:UNIT:100 block) ---\nconst paragraphs3 = body.paragraphs;\nparagraphs3.load(\"items/text\");\nawait context.sync();\n\nconst items3 = paragraphs3.items;\n// Find the run of 5 consecutive paragraphs with the RISK texts that is\n// immediately followed by \"PUMP:UNIT:100\" (that's the original, now-stale\n// copy we need to remove; the new copy we just inserted sits elsewhere).\nlet oldRiskStart = -1;\nfor (let i = 0; i + riskTexts.length < items3.length; i++) {\n  let match = true;\n  for (let j = 0; j < riskTexts.length; j++) {\n    if (items3[i + j].text !== riskTexts[j]) {\n      match = false;\n      break;\n    }\n  }\n  if (match && items3[i + riskTexts.length].text === \"PUMP:UNIT:100\") {\n    oldRiskStart = i;\n    break;\n  }\n}\nif (oldRiskStart === -1) {\n  throw new Error(\"could not locate the stale PUMP:RISK block to delete\");\n}\nfor (let j = 0; j < riskTexts.length; j++) {\n  items3[oldRiskStart + j].delete();\n}\nawait context.sync();\n\n// --- Step 4: delete the original trailing HRS..ACE block (8 paragraphs,\n//     the ones WITHOUT the two newly inserted lines) that now sits right\n//     after \"PUMP:UNIT:220\" at the end of the document. ---\nconst paragraphs4 = body.paragraphs;\nparagraphs4.load(\"items/text\");\nawait context.sync();\n\nconst items4 = paragraphs4.items;\nconst oldTailTexts = [\n  \"PUMP:HRS:103\",\n  \"PUMP:PRS:103\",\n  \"ACE:SRS:110\",\n  \"ACE:SRS:120\",\n  \"PUMP:TBV:1\",\n  \"PUMP:PRS:6\",\n  \"PUMP:DER:2\",\n  \"ACE:SRS:1000\",\n];\n\nlet oldTailStart = -1;\nfor (let i = 0; i + oldTailTexts.length <= items4.length; i++) {\n  let match = true;\n  for (let j = 0; j < oldTailTexts.length; j++) {\n    if (items4[i + j].text !== oldTailTexts[j]) {\n      match = false;\n      break;\n    }\n  }\n  if (match) {\n    oldTailStart = i;\n    break;\n  }\n}\nif (oldTailStart === -1) {\n  throw new Error(\"could not locate the stale tail block to delete\");\n}\nfor (let j = 0; j < oldTailTexts.length; j++) {\n  items4[oldTailStart + j].delete();\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Goal (per the diff):\n#   1. Move the \"PUMP:RISK:10..50\" block (5 paragraphs) from just before the\n#      \"PUMP:UNIT:100\" block to just after the intro paragraph (i.e. right\n#      before \"PUMP:URS:1\").\n#   2. Move the trailing block of 8 paragraphs\n#        PUMP:HRS:103, PUMP:PRS:103, ACE:SRS:110, ACE:SRS:120,\n#        PUMP:TBV:1, PUMP:PRS:6, PUMP:DER:2, ACE:SRS:1000\n#      (currently at the very end, after \"PUMP:UNIT:220\") to just after\n#      \"PUMP:URS:4000\" (i.e. right before the \"PUMP:UNIT:100\" block), and\n#      while doing so, insert two brand-new lines into that block:\n#        - \"PUMP:TBV:1111\" right after \"PUMP:HRS:103\"\n#        - \"PUMP:TBD:1\" right after \"PUMP:PRS:6\" (before \"PUMP:DER:2\")\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphIndexByText {\n    param($doc, [string]$text)\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs.Item($i).Range.Text -eq ($text + \"`r\")) {\n            return $i\n        }\n    }\n    throw \"Paragraph not found: $text\"\n}\n\n# Find the index of the first paragraph (starting the search at/after\n# $fromIndex) whose text, together with the following paragraphs, matches\n# the given list of texts exactly, in order.\nfunction Get-ParagraphRunIndex {\n    param($doc, [string[]]$texts, [int]$fromIndex = 1)\n    $n = $texts.Count\n    $last = $doc.Paragraphs.Count - $n + 1\n    for ($i = $fromIndex; $i -le $last; $i++) {\n        $match = $true\n        for ($j = 0; $j -lt $n; $j++) {\n            if ($doc.Paragraphs.Item($i + $j).Range.Text -ne ($texts[$j] + \"`r\")) {\n                $match = $false\n                break\n            }\n        }\n        if ($match) { return $i }\n    }\n    throw \"Paragraph run not found starting from index $fromIndex\"\n}\n\n# --- Step 1: insert the new RISK block right before \"PUMP:URS:1\" ---\n# Inserting text *before* the plain/non-bold \"PUMP:URS:1 \" paragraph (rather\n# than appending after the bold intro paragraph) keeps the new runs from\n# inheriting the intro's bold run formatting, matching the diff where these\n# new paragraphs carry no bold run properties.\n$riskTexts = @(\n    \"PUMP:RISK:10 \",\n    \"PUMP:RISK:20 \",\n    \"PUMP:RISK:30 \",\n    \"PUMP:RISK:40 \",\n    \"PUMP:RISK:50 \"\n)\n\n$urs1Index = Get-ParagraphIndexByText $d \"PUMP:URS:1 \"\n$urs1Para = $d.Paragraphs.Item($urs1Index)\n$riskBlock = ($riskTexts -join \"`r\") + \"`r\"\n$urs1Para.Range.InsertBefore($riskBlock)\n\n# --- Step 2: insert the moved/expanded HRS..ACE block right after\n#     \"PUMP:URS:4000\" (i.e. right before the paragraph that now follows it,\n#     which -- after step 1 -- is the newly-inserted \"PUMP:RISK:10 \"). ---\n$movedBlockTexts = @(\n    \"PUMP:HRS:103\",\n    \"PUMP:TBV:1111\",   # new line\n    \"PUMP:PRS:103\",\n    \"ACE:SRS:110\",\n    \"ACE:SRS:120\",\n    \"PUMP:TBV:1\",\n    \"PUMP:PRS:6\",\n    \"PUMP:TBD:1\",      # new line\n    \"PUMP:DER:2\",\n    \"ACE:SRS:1000\"\n)\n\n$urs4000Index = Get-ParagraphIndexByText $d \"PUMP:URS:4000 \"\n$afterUrs4000Para = $d.Paragraphs.Item($urs4000Index + 1)\n$movedBlock = ($movedBlockTexts -join \"`r\") + \"`r\"\n$afterUrs4000Para.Range.InsertBefore($movedBlock)\n\n# --- Step 3: delete the original RISK block (the one that was right before\n#     the \"PUMP:UNIT:100\" block). We locate it by finding the run of 5\n#     consecutive RISK paragraphs that is immediately followed by\n#     \"PUMP:UNIT:100\" -- that's the stale copy; the new copy sits elsewhere\n#     (right before \"PUMP:URS:1\") and won't match this \"followed by\" check. ---\n$oldRiskStart = 1\nwhile ($true) {\n    $oldRiskStart = Get-ParagraphRunIndex $d $riskTexts $oldRiskStart\n    $followingIndex = $oldRiskStart + $riskTexts.Count\n    $followingText = $d.Paragraphs.Item($followingIndex).Range.Text\n    if ($followingText -eq (\"PUMP:UNIT:100\" + \"`r\")) {\n        break\n    }\n    $oldRiskStart = $oldRiskStart + 1\n}\n$oldRiskStartPara = $d.Paragraphs.Item($oldRiskStart)\n$oldRiskEndPara = $d.Paragraphs.Item($oldRiskStart + $riskTexts.Count - 1)\n$oldRiskRange = $d.Range($oldRiskStartPara.Range.Start, $oldRiskEndPara.Range.End)\n$oldRiskRange.Delete()\n\n# --- Step 4: delete the original trailing HRS..ACE block (8 paragraphs,\n#     the ones WITHOUT the two newly inserted lines) that now sits right\n#     after \"PUMP:UNIT:220\" at the end of the document. ---\n$oldTailTexts = @(\n    \"PUMP:HRS:103\",\n    \"PUMP:PRS:103\",\n    \"ACE:SRS:110\",\n    \"ACE:SRS:120\",\n    \"PUMP:TBV:1\",\n    \"PUMP:PRS:6\",\n    \"PUMP:DER:2\",\n    \"ACE:SRS:1000\"\n)\n$oldTailStart = Get-ParagraphRunIndex $d $oldTailTexts 1\n$oldTailStartPara = $d.Paragraphs.Item($oldTailStart)\n$oldTailEndPara = $d.Paragraphs.Item($oldTailStart + $oldTailTexts.Count - 1)\n$oldTailRange = $d.Range($oldTailStartPara.Range.Start, $oldTailEndPara.Range.End)\n$oldTailRange.Delete()\n"}
